# Apply the "break out stock.yaml completed" edit to the DABUR.NS 1mo sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New column R ("backup") with header, styled like the other header cells.
# ---------------------------------------------------------------------------
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)
$ws.Range("R1").Value = "backup"

# ---------------------------------------------------------------------------
# 2) Fill R2:R265 with 0 (bulk array write for speed).
# ---------------------------------------------------------------------------
$n = 264
$zeros = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) { $zeros[$i,0] = 0 }
$ws.Range("R2:R265").Value = $zeros

# ---------------------------------------------------------------------------
# 3) A handful of rows keep their prior detect_structure value as the backup
#    (rows whose two_line_structure flag is set).
# ---------------------------------------------------------------------------
$ws.Range("R174").Value = 2
$ws.Range("R199").Value = 2
$ws.Range("R258").Value = 2
$ws.Range("R261").Value = 1

# ---------------------------------------------------------------------------
# 4) detect_structure (Q) gets reset to 0 on the rows that were "broken out".
# ---------------------------------------------------------------------------
$ws.Range("Q24").Value = 0
$ws.Range("Q28").Value = 0
$ws.Range("Q29").Value = 0
$ws.Range("Q38").Value = 0
$ws.Range("Q52").Value = 0

# ---------------------------------------------------------------------------
# 5) isPivot flag recomputed for row 263.
# ---------------------------------------------------------------------------
$ws.Range("O263").Value = 2

# ---------------------------------------------------------------------------
# 6) Three new months of data appended at the bottom (rows 266-268).
# ---------------------------------------------------------------------------
$ws.Range("A265").Copy()
$ws.Range("A266:A268").PasteSpecial(-4122)

$newRows = @(
  @(45474, 600.7000122070312, 662.3499755859375, 600.7000122070312, 635.7000122070312, 632.9890747070312, 65002992, 2024, 7, 1, 0, 0, 0, 27, 0, 0, 0),
  @(45505, 638.9000244140625, 655.8499755859375, 598.5999755859375, 637.1500244140625, 637.1500244140625, 53289017, 2024, 8, 1, 0, 0, 0, 31, 0, 1, 1),
  @(45536, 637.7999877929688, 672,               620.5999755859375, 633,               633,               53425949, 2024, 9, 1, 0, 0, 0, 35, 0, 0, 0)
)

$startRow = 266
for ($r = 0; $r -lt $newRows.Count; $r++) {
  $rowVals = $newRows[$r]
  $rowNum = $startRow + $r
  $rowArr = New-Object 'object[,]' 1,17
  for ($c = 0; $c -lt 17; $c++) { $rowArr[0,$c] = $rowVals[$c] }
  $ws.Range($ws.Cells.Item($rowNum,1), $ws.Cells.Item($rowNum,17)).Value = $rowArr
}
